$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - row 3: "第九届环形宇宙动漫游戏嘉年华", row 5: "心动恋章·冬日序国乙&代号鸢同人only"
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 2858
$wsExhibit.Range("F5").Value = 16

# Sheet "全部类型" (All types) - row 7 and row 10 mirror the same two events
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 2858
$wsAll.Range("F10").Value = 16
